$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the existing row 4 (CAN: section), pushing
# CAN:/1..7 down by one row (rows 4-11 -> 5-12).
$ws.Rows.Item(4).Insert()

# Add the new DIO section, leaving row 13 blank as a separator.
$ws.Range("A14").Value = "DIO:"

$ws.Range("A15").Value = 0
$ws.Range("B15").Value = "PhotoEye for Gear"

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "ball sensor"

$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "ball indicator"

$ws.Range("B17").Select()
